$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($range, $values) {
    $arr = New-Object 'object[,]' 1,$values.Count
    for ($i = 0; $i -lt $values.Count; $i++) { $arr[0,$i] = $values[$i] }
    $ws.Range($range).Value = $arr
}

function Set-ColumnValues($range, $values) {
    $arr = New-Object 'object[,]' $values.Count,1
    for ($i = 0; $i -lt $values.Count; $i++) { $arr[$i,0] = $values[$i] }
    $ws.Range($range).Value = $arr
}

# --- New daily-rate column E for rows 117-122 ---
Set-ColumnValues "E117:E122" @(0.11989087301587, 0.1210729746444, 0.11715888278388, 0.1140380713689, 0.15305435305435, 0.13011363636364)

# --- Fill data rows 211-216 (B:D and F:BE; column E intentionally stays blank, matching this sheet block's existing pattern) ---
Set-RowValues "B211:D211" @(0.091160503022997, 0.14436631900024, 0.13339322109074)
Set-RowValues "F211:BE211" @(0.076195584316284, 0.079655343431943, 0.09014129873695199, 0.07448583955876401, 0.055435220688306, 0.072669419347552, 0.085631755911595, 0.10057078838248, 0.074342342324696, 0.079998866293551, 0.12620138017397, 0.091820372021712, 0.090557550574135, 0.10046082487603, 0.11770119834686, 0.092837786185005, 0.10024097294727, 0.074180616050156, 0.06859331304727399, 0.08378107244272801, 0.07640825037325601, 0.083058633495533, 0.11171896693013, 0.080958537659389, 0.13406335881384, 0.10404795648269, 0.09499580871857601, 0.11766561364629, 0.13098338384321, 0.08032110160457499, 0.07111936328985601, 0.081329656390967, 0.072985260062744, 0.06938511668493499, 0.089518383531228, 0.13271075194164, 0.08017026410255, 0.07833361304380999, 0.093474602193795, 0.077398519806433, 0.10431931920196, 0.1157254291327, 0.10851662129536, 0.1015695130062, 0.10237736269644, 0.08026388911283899, 0.08438236306249799, 0.078667241037551, 0.07425457657885701, 0.09016508449786199, 0.087062263338543, 0.12200014425304)
Set-RowValues "B212:D212" @(0.090951197386741, 0.08758286120105099, 0.062972067285486)
Set-RowValues "F212:BE212" @(0.07041045652406, 0.065406261182633, 0.076029164488341, 0.053305859523713, 0.057072725191514, 0.05025142372602, 0.060464900833542, 0.061134088780444, 0.078113551412702, 0.058203018989916, 0.08194345058265699, 0.081017148120732, 0.06612537465036999, 0.069776258000011, 0.067738229481469, 0.06315890295990401, 0.072873671189992, 0.060806289260301, 0.055684454400932, 0.065921316340589, 0.061284829809193, 0.063166379407987, 0.06608548506715201, 0.09257585688965, 0.068382585162341, 0.08500204559993101, 0.060089849847192, 0.06742910805406099, 0.07147337611197301, 0.056325042414266, 0.053512799707223, 0.072649554829504, 0.066000410124539, 0.056674210318516, 0.060721426530042, 0.06614408514509799, 0.067952861541428, 0.05761380441479, 0.078365516118098, 0.059780145915509, 0.066794749742853, 0.075138997330714, 0.06497745595587601, 0.063773778158527, 0.085816785680188, 0.057016068681163, 0.077919078148626, 0.064957536278029, 0.070608385624934, 0.06758859925994599, 0.057423908629954, 0.09092787968611001)
Set-RowValues "B213:D213" @(0.07722624407013801, 0.067096197876128, 0.059141264201759)
Set-RowValues "F213:BE213" @(0.067077766015121, 0.060404390416145, 0.06945255468149, 0.055854996843605, 0.056943610039147, 0.056038773538383, 0.045288276575819, 0.043992166265906, 0.07307412363780499, 0.042358668509701, 0.053717656093572, 0.056438061336859, 0.045243461177322, 0.04927331248094, 0.049759111914003, 0.05050144465755, 0.056654069903717, 0.053407544393188, 0.047112859460228, 0.054447320161754, 0.047451672553534, 0.046949716619034, 0.054513180395716, 0.10523050045389, 0.052748467118161, 0.06922123530356999, 0.056698457377588, 0.049359659957464, 0.05518243538945, 0.058959137951034, 0.0609470348846, 0.073160574722023, 0.058551529823119, 0.06542034278386299, 0.05343725455796, 0.052123453206227, 0.052937494651485, 0.04791699922931, 0.077067789750097, 0.051873171191007, 0.049215770417817, 0.049412890583126, 0.049276502889787, 0.043994961224885, 0.059352313894281, 0.048918282160112, 0.071534246827664, 0.051746210459417, 0.054785301583838, 0.053578265765766, 0.051993929718068, 0.06964790790438)
Set-RowValues "B214:D214" @(0.10971337650102, 0.12197129041411, 0.13329526129196)
Set-RowValues "F214:BE214" @(0.09124387464455599, 0.083448957641741, 0.095053007095668, 0.075845440622181, 0.054481106411887, 0.075766568996347, 0.082280218300047, 0.097933890305533, 0.075901522605289, 0.08494218147575799, 0.11657017962912, 0.10608734384083, 0.087088339936137, 0.1028421199264, 0.10611414069857, 0.088201639025156, 0.099828748559909, 0.06690391354160199, 0.061393513000419, 0.076580877342703, 0.077823958284102, 0.0757675957977, 0.099704015931282, 0.079376984850653, 0.12562761528984, 0.11238346019661, 0.09116962497084, 0.11607428417781, 0.11883847008221, 0.082828474501538, 0.069434308235526, 0.088597998550419, 0.078641135021985, 0.068249641017639, 0.09183895434273701, 0.11815440814549, 0.083816888697006, 0.077039187782838, 0.081645993686478, 0.07122449322923601, 0.10115617189001, 0.1149314102667, 0.10145397692363, 0.09224189877919101, 0.11368904483561, 0.074354107283576, 0.081160091804291, 0.074134139719071, 0.075417435261334, 0.080012894509355, 0.079154674129952, 0.11751041457798)
Set-RowValues "B215:D215" @(0.071656644912978, 0.11514264674101, 0.11728369557904)
Set-RowValues "F215:BE215" @(0.06671294203488599, 0.07641698298273999, 0.09864195635749801, 0.08545064425186299, 0.062674318589397, 0.084266569331998, 0.09190908516375799, 0.10746362916289, 0.085426450462603, 0.089428654717105, 0.11413662173604, 0.11423260062487, 0.087283429716633, 0.10113022871188, 0.10622508924233, 0.088967899908086, 0.10191469289892, 0.073338231396814, 0.066086905038686, 0.082394289605163, 0.08262609941381099, 0.08218796265157401, 0.10079911283249, 0.10497854001961, 0.12817192740367, 0.11306955144533, 0.091658002710411, 0.11434055816064, 0.12486600829312, 0.082224321925604, 0.072422219350134, 0.08627675172810401, 0.078928603887918, 0.067464924133763, 0.090831133200206, 0.11674737341344, 0.088895925546859, 0.081249257507226, 0.084111620491903, 0.080085597547358, 0.1052701225728, 0.12526187971202, 0.10895378554662, 0.09621873645010399, 0.11633192437292, 0.077936794517037, 0.077848591165569, 0.08329506617303099, 0.077177500722329, 0.087529890926791, 0.08365249373913, 0.12036639729355)
Set-RowValues "B216:D216" @(0.094235027961606, 0.11839466954152, 0.12530298582001)
Set-RowValues "F216:BE216" @(0.07890550873268901, 0.072313362018451, 0.086184001311187, 0.066871109037553, 0.049184931139857, 0.070482032253063, 0.083297133424792, 0.10185374397948, 0.071895327627782, 0.084656653981367, 0.12444245829502, 0.10487086502707, 0.08603248499244499, 0.10320089527669, 0.11753072001344, 0.095136744392518, 0.11093350521934, 0.079714210708069, 0.07381781225228499, 0.091963997419264, 0.088541014062524, 0.08523826654595, 0.10608397606774, 0.085366670194547, 0.13140234472487, 0.11557281718307, 0.09107946995246601, 0.11679562486638, 0.12800388267189, 0.08274485353746899, 0.07039159785902201, 0.085304565439682, 0.07775267139183099, 0.067940466707752, 0.08668843907608, 0.11864128076216, 0.08133933593116099, 0.071090875494212, 0.086581546959404, 0.06505326336792699, 0.09878999836252, 0.11417546076272, 0.10081423248645, 0.090169475986866, 0.1130274928121, 0.07708422534616, 0.076860277021783, 0.07209607933696, 0.076801634002246, 0.092557461443025, 0.08507831424894501, 0.1259786063693)

# --- New row labels (dates) A216:A221 ---
Set-ColumnValues "A216:A221" @("02 09 2020", "03 09 2020", "04 09 2020", "05 09 2020", "06 09 2020", "07 09 2020")
